$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 08:21"

# --- Refresh COVID-19 stats for several countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4568375
$ws.Range("C4").Value = 338
$ws.Range("D4").Value = 2245521
$ws.Range("E4").Value = 2169009
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 153845

# Row 6: India
$ws.Range("B6").Value = 1587982
$ws.Range("C6").Value = 3598
$ws.Range("D6").Value = 1022565
$ws.Range("E6").Value = 530382
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 35035

# Row 36: Israel
$ws.Range("B36").Value = 68556
$ws.Range("C36").Value = 257
$ws.Range("D36").Value = 35513
$ws.Range("E36").Value = 32552

# Row 64: Uzbekistan
$ws.Range("B64").Value = 22872
$ws.Range("C64").Value = 287
$ws.Range("E64").Value = 9803
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 132

# Row 74: El Salvador
$ws.Range("D74").Value = 8095
$ws.Range("E74").Value = 7316

# Row 143: Georgia
$ws.Range("B143").Value = 1160
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 935
$ws.Range("E143").Value = 208

# Row 168: Birmania
$ws.Range("B168").Value = 353
$ws.Range("C168").Value = 2
$ws.Range("E168").Value = 53

# Row 176: Camboya
$ws.Range("B176").Value = 234
$ws.Range("C176").Value = 8
$ws.Range("D176").Value = 162
$ws.Range("E176").Value = 72

# --- San Martin (Parte Holandesa) overtakes Monaco and Aruba in the ranking ---
# New order becomes: ... Brunei, San Martin (Parte Holandesa), Monaco, Aruba, Seychelles ...
# Row 183 now holds San Martin with refreshed figures.
$ws.Range("A183").Value = "San Martin (Parte Holandesa)"
$ws.Range("B183").Value = 126
$ws.Range("C183").Value = 11
$ws.Range("D183").Value = 64
$ws.Range("E183").Value = 47
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 15

# Row 184 now holds Monaco (values unchanged, shifted down one row).
$ws.Range("A184").Value = "Monaco"
$ws.Range("B184").Value = 120
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 105
$ws.Range("E184").Value = 11
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 4

# Row 185 now holds Aruba (values unchanged, shifted down one row).
$ws.Range("A185").Value = "Aruba"
$ws.Range("B185").Value = 119
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 102
$ws.Range("E185").Value = 14
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 3
